# "Generate Report for Handoff"
#
# Refresh the localization-status report: swap the old handoff-commit
# UUID/filenames for the new ones and bump the handoff timestamps, on all
# three sheets (Overview, zh-cn, de-de). Hyperlink display text is kept in
# sync with the new cell text; hyperlink targets (Address) are left as-is,
# matching the source diff (the .rels files are untouched there too).

$wb = $excel.ActiveWorkbook

$oldId = "0242525c-20e3-4ef6-8e00-69ef082578cb"
$newId = "b427a45a-5548-4ce8-ab20-ce965e8641cf"

$oldZhHash = "482319fbc909afa334a7d4ca4278df11fd59fe00"
$newZhHash = "ac964c4e2ba753f29acaf7f317cf4ed242a0e8b3"

$oldDeHash = "482319fbc909afa334a7d4ca4278df11fd59fe00"
$newDeHash = "ac964c4e2ba753f29acaf7f317cf4ed242a0e8b3"

$oldMdName = "$oldId.md"
$newMdName = "$newId.md"

$oldZhXlf = "$oldId.$oldZhHash.zh-cn.xlf"
$newZhXlf = "$newId.$newZhHash.zh-cn.xlf"

$oldDeXlf = "$oldId.$oldDeHash.de-de.xlf"
$newDeXlf = "$newId.$newDeHash.de-de.xlf"

$oldHandoffDate = "2016-46-20 10:46:40"
$newHandoffDate = "2016-47-20 10:47:04"

$oldZhDatetime = "2016-03-20 10:46:36"
$newZhDatetime = "2016-03-20 10:47:00"

$oldDeDatetime = "2016-03-20 10:46:40"
$newDeDatetime = "2016-03-20 10:47:04"

# Original (unchanged) hyperlink targets, kept verbatim from the workbook's
# existing relationships so the link addresses stay exactly as they were.
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/5a3aa8759330847071af3ceaae9ad186abc15b7e/e2e/$oldMdName"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1cc37db4aae5ce45b6dc1a33384309d7bcb4c139/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11637ce70b53c51bfa47d199d0ef02fc0f9c4a54/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newHandoffDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName)

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhDatetime

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $mdAddress, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfAddress, "", "", $newZhXlf)

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newDeDatetime

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $mdAddress, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfAddress, "", "", $newDeXlf)

Write-Host "Report refreshed for handoff: $newMdName"
